$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

# Row 1 (table row 1, all 5 cells simple 1:1 replacements)
Replace-Text "42÷7=" "97÷3="
Replace-Text "75÷8=" "18÷3="
Replace-Text "10÷2=" "45÷9="
Replace-Text "62÷2=" "53÷6="
Replace-Text "39÷3=" "43÷8="

# Row 5 (table row 5) - cells shift: a cell is inserted after the 3rd cell and
# the last cell is removed, so we rewrite every cell in the row positionally
# rather than doing pure text substitution.
$t = $d.Tables.Item(1)
$r5 = $t.Rows.Item(5)
$r5.Cells.Item(1).Range.Text = "58÷3="
$r5.Cells.Item(2).Range.Text = "14÷6="
$r5.Cells.Item(3).Range.Text = "45÷4="
$r5.Cells.Item(4).Range.Text = "76÷4="
$r5.Cells.Item(5).Range.Text = "77÷9="

# Row 9 (table row 9, all 5 cells simple 1:1 replacements)
Replace-Text "82÷2=" "38÷6="
Replace-Text "57÷2=" "64÷3="
Replace-Text "69÷7=" "54÷7="
Replace-Text "39÷9=" "34÷5="
Replace-Text "52÷9=" "34÷4="

# Row 13 (table row 13, all 5 cells simple 1:1 replacements)
Replace-Text "30÷3=" "53÷3="
Replace-Text "58÷2=" "21÷7="
Replace-Text "91÷9=" "72÷4="
Replace-Text "92÷6=" "57÷8="
Replace-Text "18÷4=" "16÷9="

# Row 17 (table row 17, all 5 cells simple 1:1 replacements)
Replace-Text "25÷4=" "79÷9="
Replace-Text "33÷4=" "94÷5="
Replace-Text "23÷9=" "34÷4="
Replace-Text "60÷9=" "31÷3="
Replace-Text "39÷5=" "66÷3="
